$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing all existing
# episode rows down by one, just like a new RSS item being prepended to
# the "played" log while everything else keeps its relative order.
$ws.Rows.Item(2).Insert()

# The inserted row picked up row 1's (header) formatting by default;
# fix it up so column A matches the plain data-row style used
# throughout the rest of the table (copy format from the row below,
# which is the old row 2 that got pushed down to row 3), and strip the
# inherited header formatting from the other inserted cells.
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4122)
$ws.Range("B2:E2").ClearFormats()
$excel.CutCopyMode = 0

# Populate the newly inserted row 2 with the latest played episode.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "Tue, 7 Dec 2021 10:50:00 +0000"
$ws.Cells.Item(2, 3).Value = "A New Strategy for Prosecuting School Shootings"
$ws.Cells.Item(2, 4).Value = "00:23:29"
$ws.Cells.Item(2, 5).Value = "https://dts.podtrac.com/redirect.mp3/chrt.fm/track/8DB4DB/pdst.fm/e/nyt.simplecastaudio.com/03d8b493-87fc-4bd1-931f-8a8e9b945d8a/episodes/e7d47081-01f3-426c-8f22-c117171ba5f5/audio/128/default.mp3?aid=rss_feed&awCollectionId=03d8b493-87fc-4bd1-931f-8a8e9b945d8a&awEpisodeId=e7d47081-01f3-426c-8f22-c117171ba5f5&feed=54nAGcIl"

# Renumber the episode index column for all the rows that were shifted
# down (previously A2:A11 held 0..9, now living at A3:A12).
$ws.Cells.Item(3, 1).Value = 1392
$ws.Cells.Item(4, 1).Value = 1393
$ws.Cells.Item(5, 1).Value = 1394
$ws.Cells.Item(6, 1).Value = 1395
$ws.Cells.Item(7, 1).Value = 1396
$ws.Cells.Item(8, 1).Value = 1397
$ws.Cells.Item(9, 1).Value = 1398
$ws.Cells.Item(10, 1).Value = 1399
$ws.Cells.Item(11, 1).Value = 1400
$ws.Cells.Item(12, 1).Value = 1401

# Append four more historical episodes below the existing data, copying
# the column-A data style down onto the new rows first.
$ws.Range("A12").Copy()
$ws.Range("A13:A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @(1402, "Tue, 30 Nov 2021 18:05:41 GMT", "Babbage: Omicron and on", "00:28:32", "https://sphinx.acast.com/theeconomistallaudio/theeconomistbabbage/babbage-omicronandon/media.mp3"),
    @(1403, "Tue, 30 Nov 2021 10:57:53 GMT", "Centrifugal forces: Iran nuclear talks resume", "00:22:51", "https://sphinx.acast.com/theeconomistallaudio/theintelligencepodcast/centrifugalforces-irannucleartalksresume/media.mp3"),
    @(1404, "Mon, 29 Nov 2021 16:30:58 GMT", "The World Ahead: The eagle and the dragon", "00:22:53", "https://sphinx.acast.com/theeconomistallaudio/theworldahead/theworldahead-theeagleandthedragon/media.mp3")
)

$row = 13
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
